# Refresh cryptos list values (prices / 1h volume %) per upstream data pull,
# including a TRON/Chainlink row swap (rows 12 and 13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.123.16"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "2.380.50"
$ws.Range("E3").Value = "  +3.83%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'303.24"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'97.48"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("D7").Value = "'0.509"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.502"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "'34.15"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'18.58"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.122"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").Value = "'6.80"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "2.750.23"
$ws.Range("E15").Value = "  +3.87%  "
$ws.Range("D16").Value = "2.384.93"
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("E17").Value = "  +3.83%  "
$ws.Range("D18").Value = "43.115.12"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "'12.23"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'6.42"
$ws.Range("E20").Value = "  +7.34%  "
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "'236.04"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'24.80"
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'9.14"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").Value = "'31.61"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").Value = "'0.0735"
$ws.Range("E33").Value = "  +5.30%  "
$ws.Range("D34").Value = "'17.24"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  +6.69%  "
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("D40").Value = "'22.33"
$ws.Range("E40").Value = "  +11.58%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "'107.85"
$ws.Range("E42").Value = "  -34.92%  "
$ws.Range("D43").Value = "1.952.94"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").Value = "'2.13"
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  -11.47%  "
$ws.Range("D48").Value = "2.607.13"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("D49").Value = "'52.91"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("D51").Value = "'72.20"
$ws.Range("E51").Value = "  +1.64%  "
